$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in row 2
$ws.Range("A2").Value = -1.1862168312072754
$ws.Range("C2").Value = 3.819366455078125

# Update text values (shared strings) in row 2
$ws.Range("B2").Value = "Disappointment"
$ws.Range("D2").Value = "BecomeRich"
